# Sprint III test-case sheet update:
# Mark the "header" row of every multi-row test case in "Casos de Prueba"
# (column H, the "Estado" column) as "Test Aprobado", matching the pattern
# already used for the single-row test cases in rows 3-23.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Prueba")

$estado = "Test Aprobado"

$filaCabeceras = @(24, 26, 28, 30, 32, 34, 36, 39, 42, 45, 46, 50, 54)

foreach ($fila in $filaCabeceras) {
    $ws.Range("H$fila").Value = $estado
}

# Reflect where the author ended up working (bottom of the table).
[void]$ws.Activate()
[void]$ws.Range("H54:H58").Select()

Write-Output "Updated $($filaCabeceras.Count) Estado cells in 'Casos de Prueba'."
